$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Label" header in column H, matching the bold/border/center style of the
# other header cells (copy format from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Updated prediction/error values (higher precision refit results)
$ws.Range("D2").Value = 0.5563318197835855
$ws.Range("E2").Value = 0.5563318197835855

$ws.Range("D8").Value = 0.8120576269257854
$ws.Range("E8").Value = 0.1879423730742146

$ws.Range("D10").Value = 0.8330774843677082
$ws.Range("E10").Value = 0.1669225156322918

$ws.Range("D11").Value = 0.6073302221207438
$ws.Range("E11").Value = 0.3926697778792562
$ws.Range("F11").Value = 0.5403760671615601

# Label column: 0 for Control rows, 1 for MDD rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
